$d = $word.ActiveDocument

# Target the first paragraph, which contains the **ID__...__ID** bookmark text.
$p = $d.Paragraphs.Item(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space, no line style,
# matching the target <w:pBdr> with only w:space attributes.
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p.Format.LeftIndent = 11.25

# Merge the two runs ("**ID__AFFARS_mp_5306_502_topic_5__ID**" + trailing space)
# into a single run with updated, uppercased text and no trailing space.
$d.Content.Find.Execute("**ID__AFFARS_mp_5306_502_topic_5__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_MP_5306_502_4__ID**", 2)
